$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report volume number + week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  5/1/2023  Through  5/7/2023"

# --- Row 14 (Murder) ---
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = -60

# --- Row 15 (Rape) ---
$ws.Range("F15").Value = 4
$ws.Range("I15").Value = 8
$ws.Range("K15").Value = 60
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -20
$ws.Range("N15").Value = 0

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -57.142857142857
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -45
$ws.Range("I16").Value = 78
$ws.Range("J16").Value = 75
$ws.Range("K16").Value = 4
$ws.Range("L16").Value = 39.285714285714
$ws.Range("M16").Value = -18.75
$ws.Range("N16").Value = -70.229007633587

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 29
$ws.Range("H17").Value = 11.538461538461
$ws.Range("I17").Value = 131
$ws.Range("J17").Value = 114
$ws.Range("K17").Value = 14.912280701754
$ws.Range("L17").Value = 37.894736842105
$ws.Range("M17").Value = 45.555555555555
$ws.Range("N17").Value = 42.391304347826

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 21
$ws.Range("H18").Value = 110
$ws.Range("I18").Value = 84
$ws.Range("J18").Value = 47
$ws.Range("K18").Value = 78.723404255319
$ws.Range("L18").Value = 40
$ws.Range("M18").Value = -20.754716981132
$ws.Range("N18").Value = -85.053380782918

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 49
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = 2.083333333333
$ws.Range("I19").Value = 209
$ws.Range("J19").Value = 240
$ws.Range("K19").Value = -12.916666666666
$ws.Range("L19").Value = 60.769230769230
$ws.Range("M19").Value = 77.118644067796
$ws.Range("N19").Value = 6.632653061224

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 24
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 57
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 216.666666666667
$ws.Range("I20").Value = 169
$ws.Range("J20").Value = 126
$ws.Range("K20").Value = 34.126984126984
$ws.Range("L20").Value = 113.924050632911
$ws.Range("M20").Value = 108.641975308642
$ws.Range("N20").Value = -74

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 51
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = 64.516129032258
$ws.Range("F21").Value = 172
$ws.Range("H21").Value = 39.837398373983
$ws.Range("I21").Value = 681
$ws.Range("J21").Value = 610
$ws.Range("K21").Value = 11.639344262295
$ws.Range("L21").Value = 58.372093023255
$ws.Range("M21").Value = 35.387673956262
$ws.Range("N21").Value = -61.633802816901

# --- Row 22 (Transit): F22 changes from numeric 1 to text "0" ---
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "0"
$ws.Range("C22").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("H22").Value = -100
$ws.Range("M22").Value = -12.5

# --- Row 23 (Housing) ---
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 12
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = -7.692307692307
$ws.Range("I23").Value = 48
$ws.Range("J23").Value = 38
$ws.Range("K23").Value = 26.315789473684
$ws.Range("L23").Value = 84.615384615384
$ws.Range("M23").Value = 118.181818181818

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -16.666666666666
$ws.Range("F24").Value = 114
$ws.Range("G24").Value = 92
$ws.Range("H24").Value = 23.913043478260
$ws.Range("I24").Value = 526
$ws.Range("J24").Value = 410
$ws.Range("K24").Value = 28.292682926829
$ws.Range("L24").Value = 51.585014409221
$ws.Range("M24").Value = 78.911564625850

# --- Row 25 (Misd. Assault) ---
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 50
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = 19.047619047619
$ws.Range("I25").Value = 197
$ws.Range("J25").Value = 167
$ws.Range("K25").Value = 17.964071856287
$ws.Range("L25").Value = 68.376068376068
$ws.Range("M25").Value = -12.444444444444

# --- Row 26 (UCR Rape*) ---
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = 300
$ws.Range("I26").Value = 13
$ws.Range("K26").Value = 18.181818181818
$ws.Range("L26").Value = 18.181818181818

# --- Row 27 (Other Sex Crimes): D27 -> text "0", E27 -> text "***.*" ---
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 19
$ws.Range("K27").Value = 5.555555555555
$ws.Range("L27").Value = 46.153846153846

# --- Row 28 (Shooting Vic.) ---
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -85.714285714285
$ws.Range("L28").Value = -16.666666666666

# --- Row 29 (Shooting Inc.) ---
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -75
$ws.Range("L29").Value = -16.666666666666

# --- Row 30 (Hate Crimes): F30 changes from numeric 2 to text "0" ---
$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "0"
$ws.Range("C22").Copy()
$ws.Range("F30").PasteSpecial(-4122)
